$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2881169905109251
$ws.Range("C2").Value = 0.3048912486333797
$ws.Range("D2").Value = 3.223369029078222
$ws.Range("E2").Value = 0.5333859586016987
$ws.Range("G2").Value = 4.349763226824225
